$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old 4x3 (nom/prenom/age) table entirely before writing the new,
# smaller 3x2 table.
$ws.Cells.Clear()

# Write the new data. The shared-string table is populated in first-write
# order, so "B" must be written before "A" to reproduce the target
# B/A/C ordering.
$ws.Range("A2").Value = "B"
$ws.Range("A1").Value = "A"
$ws.Range("A3").Value = "C"
$ws.Range("B1").Value = 1
$ws.Range("B2").Value = 2
$ws.Range("B3").Value = 1

# Match the saved selection state (active cell B3).
$ws.Range("B3").Select()
